$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.829945333333333
$ws.Range("H2").Value = 5.489835999999999
$ws.Range("I2").Value = 0.4190796720210465
$ws.Range("J2").Value = 0.4190796720210465
$ws.Range("M2").Value = 1.168007333333333
$ws.Range("N2").Value = 3.504022
$ws.Range("O2").Value = 0.1638609704511517
$ws.Range("P2").Value = 0.1638609704511517
$ws.Range("Q2").Value = 2.137389568932444
$ws.Range("R2").Value = 19.236506120392
$ws.Range("S2").Value = 0.06867080175371905
$ws.Range("T2").Value = 0.06867080175371905

# Row 3
$ws.Range("G3").Value = 1.829945333333333
$ws.Range("H3").Value = 5.489835999999999
$ws.Range("I3").Value = 0.4190796720210465
$ws.Range("J3").Value = 0.4190796720210465
$ws.Range("O3").Value = 0.5019752511630595
$ws.Range("P3").Value = 0.5019752511630595
$ws.Range("Q3").Value = 6.547725567254665
$ws.Range("R3").Value = 58.92953010529198
$ws.Range("S3").Value = 0.2103676236200974
$ws.Range("T3").Value = 0.2103676236200974

# Row 4
$ws.Range("G4").Value = 1.829945333333333
$ws.Range("H4").Value = 5.489835999999999
$ws.Range("I4").Value = 0.4190796720210465
$ws.Range("J4").Value = 0.4190796720210465
$ws.Range("M4").Value = 1.915392333333333
$ws.Range("N4").Value = 5.746177
$ws.Range("O4").Value = 0.2687123938160456
$ws.Range("P4").Value = 0.2687123938160456
$ws.Range("Q4").Value = 3.505063261885777
$ws.Range("R4").Value = 31.545569356972
$ws.Range("S4").Value = 0.1126119018684186
$ws.Range("T4").Value = 0.1126119018684187

# Row 5
$ws.Range("G5").Value = 1.829945333333333
$ws.Range("H5").Value = 5.489835999999999
$ws.Range("I5").Value = 0.4190796720210465
$ws.Range("J5").Value = 0.4190796720210465
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.46654
$ws.Range("N5").Value = 1.39962
$ws.Range("O5").Value = 0.06545138456974327
$ws.Range("P5").Value = 0.06545138456974327
$ws.Range("Q5").Value = 0.8537426958133332
$ws.Range("R5").Value = 7.683684262319998
$ws.Range("S5").Value = 0.02742934477881139
$ws.Range("T5").Value = 0.02742934477881139

# Row 6
$ws.Range("I6").Value = 0.2833335737960661
$ws.Range("J6").Value = 0.2833335737960661
$ws.Range("M6").Value = 1.168007333333333
$ws.Range("N6").Value = 3.504022
$ws.Range("O6").Value = 0.1638609704511517
$ws.Range("P6").Value = 0.1638609704511517
$ws.Range("Q6").Value = 1.445057504792667
$ws.Range("R6").Value = 13.005517543134
$ws.Range("S6").Value = 0.0464273143636164
$ws.Range("T6").Value = 0.0464273143636164

# Row 7
$ws.Range("I7").Value = 0.2833335737960661
$ws.Range("J7").Value = 0.2833335737960661
$ws.Range("O7").Value = 0.5019752511630595
$ws.Range("P7").Value = 0.5019752511630595
$ws.Range("S7").Value = 0.1422264418692075
$ws.Range("T7").Value = 0.1422264418692075

# Row 8
$ws.Range("I8").Value = 0.2833335737960661
$ws.Range("J8").Value = 0.2833335737960661
$ws.Range("M8").Value = 1.915392333333333
$ws.Range("N8").Value = 5.746177
$ws.Range("O8").Value = 0.2687123938160456
$ws.Range("P8").Value = 0.2687123938160456
$ws.Range("Q8").Value = 2.369721479407667
$ws.Range("R8").Value = 21.327493314669
$ws.Range("S8").Value = 0.0761352428631961
$ws.Range("T8").Value = 0.07613524286319612

# Row 9
$ws.Range("I9").Value = 0.2833335737960661
$ws.Range("J9").Value = 0.2833335737960661
$ws.Range("K9").Value = 2.0
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.46654
$ws.Range("N9").Value = 1.39962
$ws.Range("O9").Value = 0.06545138456974327
$ws.Range("P9").Value = 0.06545138456974327
$ws.Range("Q9").Value = 0.5772028214599999
$ws.Range("R9").Value = 5.194825393139999
$ws.Range("S9").Value = 0.01854457470004605
$ws.Range("T9").Value = 0.01854457470004605

# Row 10
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 0.1530633333333333
$ws.Range("H10").Value = 0.45919
$ws.Range("I10").Value = 0.03505335944376924
$ws.Range("J10").Value = 0.03505335944376924
$ws.Range("M10").Value = 1.168007333333333
$ws.Range("N10").Value = 3.504022
$ws.Range("O10").Value = 0.1638609704511517
$ws.Range("P10").Value = 0.1638609704511517
$ws.Range("Q10").Value = 0.1787790957977778
$ws.Range("R10").Value = 1.60901186218
$ws.Range("S10").Value = 0.005743877496029071
$ws.Range("T10").Value = 0.005743877496029071

# Row 11
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 0.1530633333333333
$ws.Range("H11").Value = 0.45919
$ws.Range("I11").Value = 0.03505335944376924
$ws.Range("J11").Value = 0.03505335944376924
$ws.Range("O11").Value = 0.5019752511630595
$ws.Range("P11").Value = 0.5019752511630595
$ws.Range("Q11").Value = 0.5476757599366666
$ws.Range("R11").Value = 4.929081839429999
$ws.Range("S11").Value = 0.01759591891089507
$ws.Range("T11").Value = 0.01759591891089507

# Row 12
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 0.1530633333333333
$ws.Range("H12").Value = 0.45919
$ws.Range("I12").Value = 0.03505335944376924
$ws.Range("J12").Value = 0.03505335944376924
$ws.Range("M12").Value = 1.915392333333333
$ws.Range("N12").Value = 5.746177
$ws.Range("O12").Value = 0.2687123938160456
$ws.Range("P12").Value = 0.2687123938160456
$ws.Range("Q12").Value = 0.2931763351811111
$ws.Range("R12").Value = 2.63858701663
$ws.Range("S12").Value = 0.00941927212742952
$ws.Range("T12").Value = 0.009419272127429521

# Row 13
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 0.1530633333333333
$ws.Range("H13").Value = 0.45919
$ws.Range("I13").Value = 0.03505335944376924
$ws.Range("J13").Value = 0.03505335944376924
$ws.Range("K13").Value = 2.0
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.46654
$ws.Range("N13").Value = 1.39962
$ws.Range("O13").Value = 0.06545138456974327
$ws.Range("P13").Value = 0.06545138456974327
$ws.Range("Q13").Value = 0.07141016753333332
$ws.Range("R13").Value = 0.6426915078
$ws.Range("S13").Value = 0.002294290909415582
$ws.Range("T13").Value = 0.002294290909415582

# Row 14
$ws.Range("G14").Value = 1.146373333333333
$ws.Range("H14").Value = 3.43912
$ws.Range("I14").Value = 0.2625333947391181
$ws.Range("J14").Value = 0.2625333947391181
$ws.Range("M14").Value = 1.168007333333333
$ws.Range("N14").Value = 3.504022
$ws.Range("O14").Value = 0.1638609704511517
$ws.Range("P14").Value = 0.1638609704511517
$ws.Range("Q14").Value = 1.338972460071111
$ws.Range("R14").Value = 12.05075214064
$ws.Range("S14").Value = 0.04301897683778719
$ws.Range("T14").Value = 0.04301897683778719

# Row 15
$ws.Range("G15").Value = 1.146373333333333
$ws.Range("H15").Value = 3.43912
$ws.Range("I15").Value = 0.2625333947391181
$ws.Range("J15").Value = 0.2625333947391181
$ws.Range("O15").Value = 0.5019752511630595
$ws.Range("P15").Value = 0.5019752511630595
$ws.Range("Q15").Value = 4.101837277626665
$ws.Range("R15").Value = 36.91653549863999
$ws.Range("S15").Value = 0.1317852667628595
$ws.Range("T15").Value = 0.1317852667628595

# Row 16
$ws.Range("G16").Value = 1.146373333333333
$ws.Range("H16").Value = 3.43912
$ws.Range("I16").Value = 0.2625333947391181
$ws.Range("J16").Value = 0.2625333947391181
$ws.Range("M16").Value = 1.915392333333333
$ws.Range("N16").Value = 5.746177
$ws.Range("O16").Value = 0.2687123938160456
$ws.Range("P16").Value = 0.2687123938160456
$ws.Range("Q16").Value = 2.195754693804444
$ws.Range("R16").Value = 19.76179224424
$ws.Range("S16").Value = 0.07054597695700125
$ws.Range("T16").Value = 0.07054597695700127

# Row 17
$ws.Range("G17").Value = 1.146373333333333
$ws.Range("H17").Value = 3.43912
$ws.Range("I17").Value = 0.2625333947391181
$ws.Range("J17").Value = 0.2625333947391181
$ws.Range("K17").Value = 2.0
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.46654
$ws.Range("N17").Value = 1.39962
$ws.Range("O17").Value = 0.06545138456974327
$ws.Range("P17").Value = 0.06545138456974327
$ws.Range("Q17").Value = 0.5348290149333332
$ws.Range("R17").Value = 4.8134611344
$ws.Range("S17").Value = 0.01718317418147023
$ws.Range("T17").Value = 0.01718317418147023
